# Apply the changes described by the diff:
# 1. Update the "Date" metadata value on the Metadata sheet
# 2. Fix casing of "exerciceProfessionnel" -> "ExerciceProfessionnel" in the
#    element ID/Path/Base Path on the Elements sheet, and drop the trailing
#    period from the "Short"/"Definition" text for that row.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

# 1. Update Date value (row 8, column B on Metadata sheet)
$wsMeta.Range("B8").Value = "2025-10-29T11:46:56+00:00"

# 2. Update the ExerciceProfessionnel row on the Elements sheet (row 6)
$wsElements.Range("A6").Value = "OrientationParticuliere.ExerciceProfessionnel"
$wsElements.Range("B6").Value = "OrientationParticuliere.ExerciceProfessionnel"
$wsElements.Range("L6").Value = "Lien vers la classe ExerciceProfessionnel"
$wsElements.Range("M6").Value = "Lien vers la classe ExerciceProfessionnel"
$wsElements.Range("AF6").Value = "SavoirFaire.ExerciceProfessionnel"
